# Update "想去人数" (number of people interested) counts on several rows
# across the 展览 (Exhibition), 演出 (Performance), 本地生活 (Local life) and
# 全部类型 (All types) sheets. Values in column F are bumped up slightly,
# matching a refreshed data scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 901
$ws1.Range("F4").Value  = 742
$ws1.Range("F9").Value  = 1113
$ws1.Range("F10").Value = 12409
$ws1.Range("F11").Value = 680
$ws1.Range("F14").Value = 51
$ws1.Range("F16").Value = 306
$ws1.Range("F17").Value = 1829
$ws1.Range("F22").Value = 117
$ws1.Range("F23").Value = 321
$ws1.Range("F24").Value = 216
$ws1.Range("F26").Value = 102
$ws1.Range("F29").Value = 194
$ws1.Range("F31").Value = 1233

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 266
$ws2.Range("F8").Value = 115

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 847

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 847
$ws4.Range("F5").Value  = 901
$ws4.Range("F6").Value  = 742
$ws4.Range("F11").Value = 1113
$ws4.Range("F12").Value = 12409
$ws4.Range("F13").Value = 266
$ws4.Range("F14").Value = 680
$ws4.Range("F17").Value = 51
$ws4.Range("F18").Value = 306
$ws4.Range("F19").Value = 1829
$ws4.Range("F25").Value = 115
$ws4.Range("F26").Value = 115
$ws4.Range("F28").Value = 117
$ws4.Range("F32").Value = 321
$ws4.Range("F34").Value = 216
$ws4.Range("F36").Value = 102
$ws4.Range("F40").Value = 194
$ws4.Range("F44").Value = 1233
